# FloresRichardWeek7.docx Week 7 edit:
# Replace the "Theoretical Analysis example" draft (with the old Public int m()
# pseudocode walkthrough) with the finished "Performance Analysis" writeup for
# SortOfSort, preserving the _GoBack bookmark and the proofing-error markers
# that Word leaves around the flagged words/phrases.

$d = $word.ActiveDocument

$frag = '<w:p><w:r><w:t>Performance Analysis</w:t></w:r></w:p><w:p><w:r><w:t>O(n</w:t></w:r><w:r><w:t>^2</w:t></w:r><w:r><w:t>), within my version of SortOfSort, there are no nested for loops, and 1 while loop.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t>The non</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t>repeating instruction is assigning the default int’s for the counters and holder variables.</w:t></w:r></w:p><w:p><w:r><w:t>5 + 0 * n</w:t></w:r></w:p><w:p><w:r><w:t>The</w:t></w:r><w:r><w:t xml:space="preserve"> repeating instruction is to go throughout the array to find the largest within a given area. Swap around the int’s and then increment counters.</w:t></w:r></w:p><w:p><w:r><w:t>5 + 7*n</w:t></w:r></w:p><w:p><w:r><w:t>This is done 4 times before being repeated by the while loop.</w:t></w:r></w:p><w:p><w:r><w:t>5 + 28 * n</w:t></w:r></w:p><w:p><w:r><w:t>SortLeft and SortRight take the same time to complete but differ in what side of the array they modify and counter they use.</w:t></w:r><w:r><w:t xml:space="preserve"> That’s why 4 * 7 * n instead of 2 * 7 * n</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">28 * n is then done repeatedly until the array is </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SortOfSorted</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> from the while loop</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">So, 5 + </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>n(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>28 * n)</w:t></w:r></w:p>'

# $d.Content excludes the trailing section-properties mark, so replacing its
# contents with InsertXML swaps out every body paragraph in one shot while
# leaving the <w:sectPr> (page size/margins) untouched.
$d.Content.InsertXML($frag)

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
